$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Marks")

# Set E4 value to 10 (mirrors D4's value of 10)
$ws.Range("E4").Value = 10

# Set E5 to the same text as D5 ("Perfect")
$ws.Range("E5").Value = $ws.Range("D5").Value2

# Update the selection to E5:E12 (matches author's recorded selection state)
$ws.Range("E5:E12").Select()
